$d = $word.ActiveDocument

# The document begins with a Table-of-Contents content control (an SDT)
# wrapping a "Table of Contents" heading paragraph and a TOC field
# paragraph. Remove that whole block.

if ($d.ContentControls.Count -gt 0) {
    $cc = $d.ContentControls.Item(1)

    # Make sure the control (and its contents) are editable.
    $cc.LockContentControl = $false
    $cc.LockContents = $false

    $start = $cc.Range.Start

    # Unwrap the content control, turning its two paragraphs into plain
    # body paragraphs so they can be removed cleanly.
    $cc.Delete()

    # Figure out where the (now unwrapped) TOC paragraphs end: the start
    # of the paragraph that follows them both.
    $p1 = $d.Paragraphs.Item(1)
    $p2 = $d.Paragraphs.Item(2)
    $end = $p2.Range.End

    # Delete the TOC heading paragraph and the TOC field paragraph
    # entirely (including their paragraph marks).
    $r = $d.Range($start, $end)
    $r.Delete()
}
